$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.632.84'
$ws.Range('E2').Value = '  -3.31%  '

$ws.Range('D3').Value = '1.952.09'
$ws.Range('E3').Value = '  -2.44%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.013'
$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.36'
$ws.Range('E5').Value = '  -2.56%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.013'
$ws.Range('E6').Value = '  +0.06%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4767'
$ws.Range('E7').Value = '  -4.70%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4046'
$ws.Range('E8').Value = '  -4.24%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.86'
$ws.Range('E9').Value = '  -0.08%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08461'
$ws.Range('E10').Value = '  -6.19%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.060'
$ws.Range('E11').Value = '  -5.11%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.21'
$ws.Range('E12').Value = '  -4.91%  '

$ws.Range('D13').Value = '1.940.20'
$ws.Range('E13').Value = '  -4.05%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.619'
$ws.Range('E14').Value = '  -5.53%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.215'
$ws.Range('E15').Value = '  -4.08%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.016'
$ws.Range('E16').Value = '  +0.23%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001074'
$ws.Range('E17').Value = '  -3.43%  '

$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.33'
$ws.Range('E18').Value = '  -5.08%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06606'
$ws.Range('E19').Value = '  -1.00%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.68'
$ws.Range('E20').Value = '  -5.30%  '

$ws.Range('E21').Value = '  +0.07%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.831'
$ws.Range('E22').Value = '  -2.16%  '

$ws.Range('D23').Value = '28.678.55'
$ws.Range('E23').Value = '  -3.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.55'
$ws.Range('E24').Value = '  -3.57%  '

$ws.Range('E25').Value = '  -0.60%  '

$ws.Range('D26').Value = '2.193.51'
$ws.Range('E26').Value = '  -2.97%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.78'
$ws.Range('E27').Value = '  -2.49%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.21'
$ws.Range('E28').Value = '  -2.37%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.957'
$ws.Range('E29').Value = '  -7.45%  '

$ws.Range('E30').Value = '  -6.24%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '124.18'
$ws.Range('E31').Value = '  -3.23%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.002'
$ws.Range('E32').Value = '  -4.78%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09588'
$ws.Range('E33').Value = '  -3.48%  '

$ws.Range('E34').Value = '  -2.98%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.440'
$ws.Range('E35').Value = '  -8.70%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.667'
$ws.Range('E36').Value = '  -3.52%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02340'
$ws.Range('E37').Value = '  -5.21%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06232'
$ws.Range('E38').Value = '  -1.90%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.269'
$ws.Range('E39').Value = '  -3.21%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.707'
$ws.Range('E40').Value = '  -6.53%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6235'
$ws.Range('E41').Value = '  -5.01%  '

$ws.Range('E42').Value = '  -5.13%  '

$ws.Range('E43').Value = '  +0.07%  '

$ws.Range('E44').Value = '  -6.40%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.351'
$ws.Range('E45').Value = '  +3.48%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5944'
$ws.Range('E46').Value = '  -6.22%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.95'
$ws.Range('E47').Value = '  -3.60%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.081'
$ws.Range('E48').Value = '  -5.36%  '

$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000338'
$ws.Range('E49').Value = '  +0.80%  '

$ws.Range('B50').Value = 'PancakeSwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.407'
$ws.Range('E50').Value = '  -2.97%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06820'
